# Applies the "Finalized Experiments with Participant Generation" edit:
# renames each task-order sheet and refreshes the randomly-generated
# stimulus-file names in column B of each sheet.

$wb = $excel.ActiveWorkbook

# --- Rename sheets (new timestamps generated for this participant run) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16502912592046068"
$wb.Worksheets.Item(2).Name = "NB_TO-16502912617715423"
$wb.Worksheets.Item(3).Name = "RS_TO-16502912617725577"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502912618352811"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650291261924708"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502912591619687.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912591782024.csv"
$ws1.Range("B4").Value = "go_stims-16502912591823313.csv"
$ws1.Range("B5").Value = "GNG_stims-16502912592036169.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16502912606399279.csv"
$ws2.Range("B3").Value = "OB-16502912605634775.csv"
$ws2.Range("B4").Value = "ZB-match_8-165029126015198.csv"
$ws2.Range("B5").Value = "TB-16502912617522125.csv"
$ws2.Range("B6").Value = "OB-16502912606123722.csv"
$ws2.Range("B7").Value = "ZB-match_6-16502912596277783.csv"
$ws2.Range("B8").Value = "ZB-match_4-16502912596718059.csv"
$ws2.Range("B9").Value = "TB-165029126069051.csv"
$ws2.Range("B10").Value = "TB-1650291260730635.csv"

# --- Sheet 3: RS_TO --- (no content changes, only sheet name updated above)

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502912617867548.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912617745419.csv"
$ws4.Range("B4").Value = "MM_stims-16502912618183393.csv"
$ws4.Range("B5").Value = "ZM_stims-1650291261787709.csv"
$ws4.Range("B6").Value = "MM_stims-16502912618342755.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912618183393.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16502912618510482.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502912618765507.csv"
$ws5.Range("B4").Value = "SAT_stims-16502912618372746.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502912619074228.csv"
